# createAccount sheet cleanup:
#  - drop the "Email V" column (D) together with its hyperlinked sample value
#  - drop all the extra demo columns after "Year" (old H:V -> First Name2 .. Email caso alterno)
#  - change the sample "Title" value from "Mujer" to "Hombre"
#  - the sheet keeps its original scroll/active state (createAccount stays the active tab)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("createAccount")

# Make this sheet active again (it was the active tab in the source file) and
# scroll it so column G is visible at the left edge, matching the saved view.
$ws.Activate()

# Drop the two sample hyperlinks (D2 -> email, V2 -> email) before the cells
# that carry them are removed/restructured.
$ws.Hyperlinks.Delete()

# Column D ("Email V" / the hyperlinked sample address) is removed entirely;
# everything to its right shifts one column to the left.
$ws.Range("D1:D2").EntireColumn.Delete()

# After the shift, the headers/values we want to keep occupy A:G
# (Title, First Name, Last Name, Password, Day, Month, Year). Remove the
# remaining extra demo columns (old "First Name2" ... "Email caso alterno").
$ws.Range("H1:U2").EntireColumn.Delete()

# The sample "Title" value changes from "Mujer" to "Hombre".
$ws.Range("A2").Value = "Hombre"

$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
